$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week 5 (rows 44-51): fill in the Raspberry Pi / documentation entries ---
# (written first so the new shared strings land in the same order as the
# authored workbook: 28 installation, 29 install doc, 30 usage doc,
# 31 ergonomics, 32 cahier des charges, 33 discussion)
$ws.Range("B44").Value = "installation Raspberry Pi"
$ws.Range("D44").Value = 2
$ws.Range("B45").Value = "documentation d'installation"
$ws.Range("D45").Value = 2

$ws.Range("B48").Value = "documentation d'utilisation"
$ws.Range("D48").Value = 1.5
$ws.Range("B49").Value = "vérifier l'ergonomie"
$ws.Range("D49").Value = 1

# --- Week 1 (rows 6-9): add the two missing "Travail effectue" entries ---
$ws.Range("B7").Value = "Lecture du cahier des charges"
$ws.Range("D7").Value = 1
$ws.Range("B8").Value = "Discussion du projet"
$ws.Range("D8").Value = 2

# Re-create the medium top border under row 6 (separator above the newly
# filled row 7), matching the thick bottom rule already used elsewhere in
# the sheet.
$ws.Range("B7:C7").Borders.Item(8).Weight = -4138

# --- sheet view bookkeeping ---
$excel.ActiveWindow.Zoom = 115
$ws.Range("F12").Select()
